# Auto-generated edit script applying crypto price/volume/coin updates
# per commit message: "Updated cryptos list on Sat Nov 11 17:46:03 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '37.134.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -0.32%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '2.076.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.81%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.07%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '252.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +1.02%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '0.674'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +2.29%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '59.10'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +14.22%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +0.01%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +4.77%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '61.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +0.14%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.0791'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +6.54%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +2.47%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '16.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +6.44%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '2.380.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -0.76%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '0.817'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -1.95%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '5.55'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +8.44%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '2.073.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.93%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '37.077.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.40%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '15.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +10.61%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '74.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +3.46%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '0.0₃0926'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +10.45%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '5.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +4.55%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '239.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -0.44%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -0.22%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '2.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -2.03%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '2.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +14.29%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '169.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -0.45%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '9.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +0.81%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '20.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -1.90%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +2.72%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'" + '1.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +6.65%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'" + '4.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +7.38%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '0.0636'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +4.30%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '4.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +8.86%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +0.67%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -0.12%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '2.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -0.19%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +30.32%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -4.93%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '1.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +1.51%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = "'" + '18.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -1.65%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'" + '0.0227'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +1.28%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').Value = "'" + '4.60'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +21.58%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '1.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +0.57%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '98.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +0.04%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D47').Value = "'" + '4.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +12.79%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '2.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +8.75%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '2.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +0.15%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '1.302.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -1.23%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '6.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -0.30%  '
$ws.Range('E51').Style = 'Normal'
